# backup_list.xlsx - add new item(s) to the workbook
#
# Adds three new worksheets ("信息收集", "弱口令生成", "值得关注的github")
# after the existing two sheets, fills them with tool/reference entries,
# and appends several new rows to the first sheet ("企业内建安全工具").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the three new worksheets, in order, at the end of the
#    workbook so sheetId/relationship order matches (3, 4, 5).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "信息收集"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "弱口令生成"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "值得关注的github"

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 2. Sheet 1 ("企业内建安全工具") - append new tool rows.
# ---------------------------------------------------------------------
$ws1.Range("B5").Value = "https://github.com/Ekultek/Zeus-Scanner.git"
$ws1.Range("A5").Value = "Zeus-Scanner"
$ws1.Range("C5").Value = "http://www.freebuf.com/sectool/158355.html"

$ws1.Range("B6").Value = "https://github.com/alpha1e0/pentestdb.git"
$ws1.Range("C6").Value = "WEB渗透测试数据库"
$ws1.Range("A6").Value = "alpha1e0/pentestdb"

$ws1.Range("A8").Value = "资产验活工具"

$ws1.Range("B9").Value = "https://bitbucket.org/LaNMaSteR53/peepingtom.git"
$ws1.Range("A9").Value = "peepingtom"
$ws1.Range("C9").Value = "资产验活/自动截图"
$ws1.Range("C10").Value = "资产验活/自动截图"
$ws1.Range("B10").Value = "https://github.com/ChrisTruncer/EyeWitness.git"
$ws1.Range("A10").Value = "EyeWitness"

$ws1.Range("A12").Value = "Github监控"

$ws1.Range("A13").Value = "GSIL"
$ws1.Range("B13").Value = "https://github.com/FeeiCN/GSIL.git"
$ws1.Range("C13").Value = "实时监控Github敏感信息泄露，并发送告警通知"

# ---------------------------------------------------------------------
# 3. Sheet 3 ("信息收集")
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = "名称"
$ws3.Range("B1").Value = "URL"
$ws3.Range("C1").Value = "备注"

$ws3.Range("A2").Value = "taoman"
$ws3.Range("C2").Value = "快速收集 https://src.edu-info.edu.cn 平台学校子域名工具"
$ws3.Range("B2").Value = "https://github.com/LandGrey/taoman.git"

$ws3.Range("A11").Value = "httpscan"
$ws3.Range("A10").Value = "轻量级信息收集工具"
$ws3.Range("B11").Value = "https://github.com/5up3rc/httpscan.git"
$ws3.Range("C11").Value = "httpscan是一个扫描指定CIDR网段的Web主机的小工具。和端口扫描器不一样，httpscan是以爬虫的方式进行Web主机发现，因此相对来说不容易被防火墙拦截"

# ---------------------------------------------------------------------
# 4. Sheet 4 ("弱口令生成")
# ---------------------------------------------------------------------
$ws4.Range("A1").Value = "名称"
$ws4.Range("B1").Value = "URL"
$ws4.Range("C1").Value = "备注"

$ws4.Range("A2").Value = "genpAss"
$ws4.Range("B2").Value = "https://github.com/Tigascan/genpAss.git"
$ws4.Range("C2").Value = "中国特色的弱口令生成器"

# ---------------------------------------------------------------------
# 5. Sheet 5 ("值得关注的github")
# ---------------------------------------------------------------------
$ws5.Range("A1").Value = "https://github.com/5up3rc"
$ws5.Range("A4").Value = "https://github.com/jas502n"
$ws5.Range("A2").Value = "https://github.com/RicterZ"

# ---------------------------------------------------------------------
# 6. Column widths (the host quantizes ColumnWidth to 1/7-character
#    steps, so these are the closest reproducible values to the
#    target pixel widths).
# ---------------------------------------------------------------------
$ws1.Columns.Item(2).ColumnWidth = 49.857142857142854
$ws1.Columns.Item(3).ColumnWidth = 40.285714285714285

$ws3.Columns.Item(1).ColumnWidth = 21.57142857142857
$ws3.Columns.Item(2).ColumnWidth = 46.42857142857143
$ws3.Columns.Item(3).ColumnWidth = 53.42857142857143

$ws4.Columns.Item(1).ColumnWidth = 29.285714285714285
$ws4.Columns.Item(2).ColumnWidth = 56.285714285714285
$ws4.Columns.Item(3).ColumnWidth = 43.14285714285714

$ws5.Columns.Item(1).ColumnWidth = 32.57142857142857

# ---------------------------------------------------------------------
# 7. Page setup for the sheets that gained a pageSetup/printer entry.
# ---------------------------------------------------------------------
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

$ws5.PageSetup.PaperSize = 9
$ws5.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 8. Selection / active cell bookkeeping, finishing back on sheet 1.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A1:C1").Select()

$ws3.Activate()
$ws3.Range("A1:C1").Select()

$ws4.Activate()
$ws4.Range("B8").Select()

$ws5.Activate()
$ws5.Range("B9").Select()

$ws1.Activate()
$ws1.Range("C13").Select()
